$wb = $excel.ActiveWorkbook

# Sheet 1: GNG_TO
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16502912089203875"
$ws1.Range("B2").Value = "go_stims-16502912088603883.csv"
$ws1.Range("B3").Value = "GNG_stims-1650291208887391.csv"
$ws1.Range("B4").Value = "go_stims-16502912088893902.csv"
$ws1.Range("B5").Value = "GNG_stims-16502912089183867.csv"

# Sheet 2: NB_TO
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16502912115600796"
$ws2.Range("B2").Value = "TB-1650291211519084.csv"
$ws2.Range("B3").Value = "OB-16502912100747805.csv"
$ws2.Range("B4").Value = "OB-16502912101467779.csv"
$ws2.Range("B5").Value = "ZB-match_0-16502912094877844.csv"
$ws2.Range("B6").Value = "OB-16502912107791076.csv"
$ws2.Range("B7").Value = "TB-16502912115450861.csv"
$ws2.Range("B8").Value = "ZB-match_6-16502912098057795.csv"
$ws2.Range("B9").Value = "ZB-match_6-16502912091503866.csv"
$ws2.Range("B10").Value = "TB-16502912113230803.csv"

# Sheet 3: RS_TO (name only)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16502912115620792"

# Sheet 4: TOL_TO
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16502912116250763"
$ws4.Range("B2").Value = "MM_stims-1650291211576082.csv"
$ws4.Range("B3").Value = "ZM_stims-16502912115630772.csv"
$ws4.Range("B4").Value = "MM_stims-16502912116080775.csv"
$ws4.Range("B5").Value = "ZM_stims-16502912115770767.csv"
$ws4.Range("B6").Value = "MM_stims-16502912116240768.csv"
$ws4.Range("B7").Value = "ZM_stims-16502912116090848.csv"

# Sheet 5: vSAT_TO
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-1650291211701076"
$ws5.Range("B2").Value = "vSAT_stims-16502912116710756.csv"
$ws5.Range("B3").Value = "SAT_stims-16502912116551144.csv"
$ws5.Range("B4").Value = "SAT_stims-16502912116290798.csv"
$ws5.Range("B5").Value = "vSAT_stims-1650291211686077.csv"
